# This workbook lists daily price/volume records for "Bruselas (repollito)"
# at "Feria Lagunitas de Puerto Montt". The edit re-associates the weekly
# observations (Fecha / Volumen / Precio minimo / Precio maximo /
# Precio promedio ponderado / Precio $/Kg) with different rows - i.e. the
# set of per-row records (columns D, J, K, L, M, P) for rows 2..27 is
# permuted. No rows are added or removed and no other columns change.
#
# $targetRow = $sourceRow means: the new contents of $targetRow (for the
# moving columns) are the OLD contents of $sourceRow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 10
    3  = 12
    4  = 15
    5  = 25
    6  = 27
    7  = 8
    8  = 19
    9  = 17
    10 = 11
    11 = 13
    12 = 7
    13 = 9
    14 = 2
    15 = 20
    16 = 14
    17 = 6
    18 = 3
    19 = 5
    20 = 21
    21 = 24
    22 = 16
    23 = 23
    24 = 22
    25 = 18
    26 = 4
    27 = 26
}

$movingCols = @("D", "J", "K", "L", "M", "P")

# First snapshot every source cell's current value, since several target
# rows read from rows that are themselves overwritten later in the loop.
$snapshot = @{}
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($col in $movingCols) {
        $key = "$col$sourceRow"
        if (-not $snapshot.ContainsKey($key)) {
            $snapshot[$key] = $ws.Range($key).Value2
        }
    }
}

# Now write the permuted values back into the sheet.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($col in $movingCols) {
        $ws.Range("$col$targetRow").Value2 = $snapshot["$col$sourceRow"]
    }
}
